$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, matching the style of the existing headers (E1)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill F2:H16 with boolean outlier flags (all FALSE except F4 which is TRUE)
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}

$ws.Cells.Item(4, 6).Value = $true
